$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 31, shifting existing rows 31-78 down to 32-79.
$ws.Rows.Item(31).Insert()

# Populate the new row 31 with the new price entry.
$ws.Range("A31").Value = 3
$ws.Range("B31").Value = "Femacal de La Calera"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44799
$ws.Range("E31").Value = 5
$ws.Range("F31").Value = 100112035
$ws.Range("G31").Value = "Bruselas (repollito)"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 55
$ws.Range("K31").Value = 15000
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = 15000
$ws.Range("N31").Value = "$/malla 15 kilos"
$ws.Range("O31").Value = "Provincia de Quillota"
$ws.Range("P31").Value = 1000
$ws.Range("Q31").Value = 15
$ws.Range("R31").Value = "Hortaliza"
